$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction6")
$ws.Activate()

# Clear existing content in row 1 from C1:O1
$ws.Range("C1:O1").Clear()

# Set new values
$ws.Range("A1").Value = 12
$ws.Range("B1").Value = 13
